# Fix damage property type: money -> diamond
# Updates columns E and F for rows 11-70 (PlayerAtt1..PlayerAtt60) following
# the new pricing curve: value = 1050 + (row - 21) * 50
# Also updates the active selection on the sheet to G30.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 11; $row -le 70; $row++) {
    $value = 1050 + ($row - 21) * 50
    $ws.Cells.Item($row, 5).Value = $value   # Column E
    $ws.Cells.Item($row, 6).Value = $value   # Column F
}

# Update the active selection shown when the sheet is opened.
$ws.Range("G30").Select()
